$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 387, pushing existing rows 387-507 down to 388-508.
$ws.Rows.Item(387).Insert()

# Populate the newly inserted row 387 with the new daily price record.
$ws.Cells.Item(387, 1).Value = 3
$ws.Cells.Item(387, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(387, 3).Value = "Coquimbo"
$ws.Cells.Item(387, 4).Value = 44985
$ws.Cells.Item(387, 5).Value = 5
$ws.Cells.Item(387, 6).Value = 100114013
$ws.Cells.Item(387, 7).Value = "Zanahoria"
$ws.Cells.Item(387, 8).Value = "Sin especificar"
$ws.Cells.Item(387, 9).Value = "Primera"
$ws.Cells.Item(387, 10).Value = 250
$ws.Cells.Item(387, 11).Value = 8500
$ws.Cells.Item(387, 12).Value = 9000
$ws.Cells.Item(387, 13).Value = 8760
$ws.Cells.Item(387, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(387, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(387, 16).Value = 438
$ws.Cells.Item(387, 17).Value = 20
$ws.Cells.Item(387, 18).Value = "Hortaliza"
